# Generate Report for Handback
# Refresh the handback-status report with the latest handoff/handback
# timestamps produced by this run (50e99dc3-... file).

$wb = $excel.ActiveWorkbook

# Overview sheet: bump the "Latest HO Xliff Generate Date" for the
# 50e99dc3 file (row 2) to reflect the new xliff generation.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-30 22:52:54"

# zh-cn sheet: the 50e99dc3 row (row 2) gets new handoff/handback times.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-30 22:52:49"
$zhcn.Range("K2").Value = "2016-08-30 22:53:15"

# de-de sheet: same for the 50e99dc3 row (row 2).
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-30 22:52:54"
$dede.Range("K2").Value = "2016-08-30 22:53:23"
